# Apply corrections to "df de autos y motos" (arrivals/offer ratios)
# across the 5 parking-lot sheets, per commit: "correccion en df de autos y motos"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('CENTENARIO')
$ws.Range("G10").Value = 1.733
$ws.Range("G11").Value = 0.1
$ws.Range("G17").Value = 0.3

$ws = $wb.Worksheets.Item('EL_PEÑON')
$ws.Range("E10").Value = 3.714
$ws.Range("F10").Value = 0.286
$ws.Range("H10").Value = 0.2
$ws.Range("I10").Value = 0.2
$ws.Range("G14").Value = 14.6
$ws.Range("F15").Value = 3.25
$ws.Range("I15").Value = 0.2
$ws.Range("D18").Value = 2
$ws.Range("G18").Value = 4.4
$ws.Range("D24").Value = 5.286
$ws.Range("F30").Value = 6.5
$ws.Range("D32").Value = 6.143
$ws.Range("G32").Value = 4.8
$ws.Range("D40").Value = 16.6
$ws.Range("E40").Value = 13.8
$ws.Range("G40").Value = 15.6
$ws.Range("H40").Value = 0.2
$ws.Range("D43").Value = 0.625
$ws.Range("G43").Value = 1.2
$ws.Range("E46").Value = 6.857
$ws.Range("H46").Value = 0.2
$ws.Range("D53").Value = 12.429

$ws = $wb.Worksheets.Item('GRANADA')
$ws.Range("D10").Value = 12.4
$ws.Range("G10").Value = 0.8
$ws.Range("I16").Value = 3.6
$ws.Range("D17").Value = 28.4
$ws.Range("G17").Value = 0.1
$ws.Range("D23").Value = 3.417
$ws.Range("G23").Value = 1.8
$ws.Range("I25").Value = 1.533
$ws.Range("D29").Value = 4.167
$ws.Range("G29").Value = 0.9
$ws.Range("H29").Value = 0.7
$ws.Range("H30").Value = 0.4
$ws.Range("D31").Value = 6.417
$ws.Range("H31").Value = 0.6
$ws.Range("I31").Value = 2.8
$ws.Range("D32").Value = 3.538
$ws.Range("G32").Value = 0.6
$ws.Range("D38").Value = 12.25
$ws.Range("G38").Value = 0.1
$ws.Range("I40").Value = 4.6
$ws.Range("G42").Value = 0.533
$ws.Range("D47").Value = 9.667
$ws.Range("G47").Value = 6.6
$ws.Range("G52").Value = 3.2
$ws.Range("D68").Value = 14.455
$ws.Range("G68").Value = 3.5
$ws.Range("H68").Value = 1.7
$ws.Range("I68").Value = 1
$ws.Range("D69").Value = 7.4
$ws.Range("G69").Value = 6.1
$ws.Range("D73").Value = 6.438
$ws.Range("G73").Value = 3
$ws.Range("D75").Value = 1.692
$ws.Range("G75").Value = 2.267
$ws.Range("I75").Value = 2.8
$ws.Range("I78").Value = 1.467
$ws.Range("D79").Value = 5
$ws.Range("G79").Value = 0.2
$ws.Range("D81").Value = 8
$ws.Range("G81").Value = 0.067
$ws.Range("D86").Value = 0.636
$ws.Range("G86").Value = 4.5
$ws.Range("I89").Value = 14.8
$ws.Range("G91").Value = 1.6
$ws.Range("D92").Value = 2

$ws = $wb.Worksheets.Item('SAN_ANTONIO')
$ws.Range("D7").Value = 6.533
$ws.Range("G7").Value = 1.667
$ws.Range("G20").Value = 3.4
$ws.Range("D30").Value = 2.857
$ws.Range("G30").Value = 2.4
$ws.Range("D58").Value = 2.375
$ws.Range("G58").Value = 2.4
$ws.Range("D59").Value = 3.615
$ws.Range("G59").Value = 0.8
$ws.Range("D60").Value = 3.833
$ws.Range("G60").Value = 3.2
$ws.Range("D78").Value = 4.5

$ws = $wb.Worksheets.Item('SAN_FERNANDO_PARQUE_DEL_PERRO')
$ws.Range("D6").Value = 3
$ws.Range("G6").Value = 1.4
$ws.Range("E8").Value = 23
$ws.Range("G18").Value = 1.2
$ws.Range("E19").Value = 0.182
$ws.Range("H19").Value = 0.1
$ws.Range("D22").Value = 11.143
$ws.Range("G22").Value = 1
$ws.Range("G28").Value = 2.8
$ws.Range("E31").Value = 2.5
$ws.Range("H31").Value = 0.2
$ws.Range("F35").Value = 0.867
$ws.Range("I35").Value = 0.067
$ws.Range("D36").Value = 3.75
$ws.Range("G36").Value = 0.1
$ws.Range("D46").Value = 1
$ws.Range("G46").Value = 1.2
$ws.Range("D48").Value = 2.8
$ws.Range("E48").Value = 2
$ws.Range("G48").Value = 2.3
$ws.Range("H48").Value = 0.1
$ws.Range("G55").Value = 1
$ws.Range("G60").Value = 4.2
$ws.Range("E65").Value = 5.143
$ws.Range("H65").Value = 0.2
$ws.Range("D67").Value = 9.75
$ws.Range("G67").Value = 1
$ws.Range("E72").Value = 6.125
$ws.Range("H72").Value = 0.2
$ws.Range("G77").Value = 0
$ws.Range("D78").Value = 4.75
$ws.Range("G78").Value = 0.8
$ws.Range("D93").Value = 4.6
$ws.Range("G93").Value = 3.6
$ws.Range("D94").Value = 0.5
$ws.Range("G94").Value = 3.3
$ws.Range("D99").Value = 8.833
$ws.Range("G99").Value = 0.2
